$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "vinayak"
$ws.Range("F2").Value = "Alice is interested in the premium package that includes advanced analytics, priority support, and additional storage capacity. She wants a detailed demo before making the decision."
$ws.Range("M2").Value = "nan"
$ws.Range("N2").Value = "nan"

# --- Add new row 3 ---
$ws.Range("A3").Value = 2

# B3, D3, H3 are numeric-looking strings that must stay text (like B2/D2/H2).
# Force text format, assign, then clear the format so no extra style sticks.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "3"
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Value = "vipul"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "9977183691"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "vinayak_sharma@technologymindz.com"
$ws.Range("F3").Value = "Charlie is evaluating enterprise-level solutions with a strong emphasis on scalability, integration with his existing ERP system, and compliance with international data protection regulations. He also needs a custom training program for his team."

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = ""
$ws.Range("G3").ClearFormats()

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "91"
$ws.Range("H3").ClearFormats()

$ws.Range("I3").Value = "Real Estate"
$ws.Range("J3").Value = "XYZ Company Ltd."
$ws.Range("K3").Value = "Berlin, Germany"
$ws.Range("L3").Value = "yes"
$ws.Range("M3").Value = "Charlie has shown strong interest in a long-term partnership if the enterprise solution aligns with his company’s compliance and integration needs. He mentioned that decision-making will involve multiple stakeholders, and the procurement cycle might take up to three months. We should prepare detailed documentation, case studies, and a tailored presentation for his board of directors."
$ws.Range("N3").Value = "nan"

Write-Output "applied edits"
